# The document's table captions follow the pattern "表 8-2-<n> <title>".
# Most of them (表 8-2-3 through 表 8-2-10) are prefixed with a stray
# "▲" marker run that the earlier captions (表 8-2-1, 表 8-2-2) do not
# have. This edit removes that leftover "▲" marker run from in front of
# the "表 8-2-3 回收紀錄" caption only, bringing it in line with the
# other "clean" captions, while leaving the caption text run itself (and
# every other "▲" marker elsewhere in the document) untouched.

$d = $word.ActiveDocument

# Anchor on the unique phrase so we grab the right occurrence (there are
# several "▲表 8-2-N ..." captions later in the document that must stay
# untouched).
$rng = $d.Content
$found = $rng.Find.Execute("▲表 8-2-3 回收紀錄", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # $rng now spans the matched phrase; isolate just the leading "▲"
    # character so the following caption run (and its own run
    # properties) is left completely intact, instead of merging the two
    # runs together via a Find/Replace.
    $markerStart = $rng.Start
    $marker = $d.Range($markerStart, $markerStart + 1)
    $marker.Delete()
}
